# Apply cell value updates per the target diff (cosinor per-day fixed-period results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.73000000000058"
$ws.Range("H2").Value = [double]"0.3764005478658895"
$ws.Range("I2").Value = [double]"0.3764005478658895"
$ws.Range("L2").Value = [double]"4.577578995159337"
$ws.Range("M2").Value = "[-4.5567050361605475, 13.711863026479222]"
$ws.Range("N2").Value = [double]"0.3182020900163158"
$ws.Range("O2").Value = [double]"0.3182020900163158"
$ws.Range("P2").Value = [double]"-0.5031579825569237"
$ws.Range("Q2").Value = "[-3.641605898755736, 2.6352899336418885]"
$ws.Range("R2").Value = [double]"0.7482648491358694"
$ws.Range("S2").Value = [double]"0.7482648491358694"
$ws.Range("T2").Value = [double]"12.1496260554843"
$ws.Range("U2").Value = "[7.526854231484894, 16.772397879483712]"
$ws.Range("V2").Value = [double]"3.442033573719527e-06"
$ws.Range("W2").Value = [double]"3.442033573719527e-06"
$ws.Range("X2").Value = [double]"2.060460460460508"
$ws.Range("Y2").Value = [double]"-10.79166166166191"
$ws.Range("Z2").Value = [double]"14.91258258258292"
$ws.Range("F3").Value = [double]"25.73000000000058"
$ws.Range("H3").Value = [double]"0.06827688251845399"
$ws.Range("I3").Value = [double]"0.06827688251845399"
$ws.Range("L3").Value = [double]"9.032108380839528"
$ws.Range("M3").Value = "[-1.2591578434518844, 19.32337460513094]"
$ws.Range("N3").Value = [double]"0.08389561980306026"
$ws.Range("O3").Value = [double]"0.08389561980306026"
$ws.Range("P3").Value = [double]"-2.012631930227695"
$ws.Range("Q3").Value = "[-3.798842768304774, -0.22642109215061534]"
$ws.Range("R3").Value = [double]"0.02808499684366628"
$ws.Range("S3").Value = [double]"0.02808499684366628"
$ws.Range("T3").Value = [double]"12.97382987083287"
$ws.Range("U3").Value = "[7.545155841381632, 18.402503900284117]"
$ws.Range("V3").Value = [double]"1.707455987998152e-05"
$ws.Range("W3").Value = [double]"1.707455987998152e-05"
$ws.Range("X3").Value = [double]"8.241841841842028"
$ws.Range("Y3").Value = [double]"0.9272072072072284"
$ws.Range("Z3").Value = [double]"15.55647647647683"
$ws.Range("F4").Value = [double]"25.73000000000058"
$ws.Range("H4").Value = [double]"0.1518944053869346"
$ws.Range("I4").Value = [double]"0.1518944053869346"
$ws.Range("L4").Value = [double]"6.568282653103847"
$ws.Range("M4").Value = "[-1.8971656949451106, 15.033731001152804]"
$ws.Range("N4").Value = [double]"0.1251223903122303"
$ws.Range("O4").Value = [double]"0.1251223903122303"
$ws.Range("P4").Value = [double]"-1.660421342437848"
$ws.Range("Q4").Value = "[-3.3648690083494284, 0.04402632347373192]"
$ws.Range("R4").Value = [double]"0.05595572814988348"
$ws.Range("S4").Value = [double]"0.05595572814988348"
$ws.Range("T4").Value = [double]"12.95638156387154"
$ws.Range("U4").Value = "[8.094528713444088, 17.81823441429899]"
$ws.Range("V4").Value = [double]"2.682908767592806e-06"
$ws.Range("W4").Value = [double]"2.682908767592806e-06"
$ws.Range("X4").Value = [double]"6.799519519519674"
$ws.Range("Y4").Value = [double]"-0.1802902902902979"
$ws.Range("Z4").Value = [double]"13.77932932932965"
$ws.Range("F5").Value = [double]"25.73000000000058"
$ws.Range("H5").Value = [double]"0.04624078557834865"
$ws.Range("I5").Value = [double]"0.04624078557834865"
$ws.Range("L5").Value = [double]"8.459006639454344"
$ws.Range("M5").Value = "[0.6789713600100669, 16.23904191889862]"
$ws.Range("N5").Value = [double]"0.03375782484802436"
$ws.Range("O5").Value = [double]"0.03375782484802436"
$ws.Range("P5").Value = [double]"-3.081842643161158"
$ws.Range("Q5").Value = "[-4.377474448245237, -1.7862108380770794]"
$ws.Range("R5").Value = [double]"1.839702291817069e-05"
$ws.Range("S5").Value = [double]"1.839702291817069e-05"
$ws.Range("T5").Value = [double]"14.82409801705776"
$ws.Range("U5").Value = "[10.284345819953591, 19.363850214161936]"
$ws.Range("V5").Value = [double]"4.332854008914921e-08"
$ws.Range("W5").Value = [double]"4.332854008914921e-08"
$ws.Range("X5").Value = [double]"12.62032032032061"
$ws.Range("Y5").Value = [double]"7.3146346346348"
$ws.Range("Z5").Value = [double]"17.92600600600641"
$ws.Range("F6").Value = [double]"25.73000000000058"
$ws.Range("H6").Value = [double]"0.4177761444968628"
$ws.Range("I6").Value = [double]"0.4177761444968628"
$ws.Range("L6").Value = [double]"5.204360665291339"
$ws.Range("M6").Value = "[-5.645062616809526, 16.053783947392205]"
$ws.Range("N6").Value = [double]"0.3391352884269776"
$ws.Range("O6").Value = [double]"0.3391352884269776"
$ws.Range("P6").Value = [double]"-2.364842518017542"
$ws.Range("Q6").Value = "[-5.497000959434392, 0.767315923399309]"
$ws.Range("R6").Value = [double]"0.1353341327823774"
$ws.Range("S6").Value = [double]"0.1353341327823774"
$ws.Range("T6").Value = [double]"17.23578479346313"
$ws.Range("U6").Value = "[11.782422779016645, 22.689146807909605]"
$ws.Range("V6").Value = [double]"8.929752381980904e-08"
$ws.Range("W6").Value = [double]"8.929752381980904e-08"
$ws.Range("X6").Value = [double]"9.684164164164383"
$ws.Range("Y6").Value = [double]"-3.142202202202274"
$ws.Range("Z6").Value = [double]"22.51053053053104"
$ws.Range("F7").Value = [double]"25.73000000000058"
$ws.Range("H7").Value = [double]"0.4000519368547175"
$ws.Range("I7").Value = [double]"0.4000519368547175"
$ws.Range("L7").Value = [double]"4.335722344331831"
$ws.Range("M7").Value = "[-4.00918833048633, 12.680633019149992]"
$ws.Range("N7").Value = [double]"0.3009384890569233"
$ws.Range("O7").Value = [double]"0.3009384890569233"
$ws.Range("P7").Value = [double]"2.748500479717197"
$ws.Range("Q7").Value = "[-0.38994743648161556, 5.886948395916009]"
$ws.Range("R7").Value = [double]"0.08454481195123598"
$ws.Range("S7").Value = [double]"0.08454481195123598"
$ws.Range("T7").Value = [double]"11.38200469029026"
$ws.Range("U7").Value = "[6.8938949609039994, 15.870114419676518]"
$ws.Range("V7").Value = [double]"6.418636561811297e-06"
$ws.Range("W7").Value = [double]"6.418636561811297e-06"
$ws.Range("X7").Value = [double]"14.47473473473506"
$ws.Range("Y7").Value = [double]"1.622612612612649"
$ws.Range("Z7").Value = [double]"27.32685685685748"
$ws.Range("F8").Value = [double]"25.73000000000058"
$ws.Range("H8").Value = [double]"0.6115400348959719"
$ws.Range("I8").Value = [double]"0.6115400348959719"
$ws.Range("L8").Value = [double]"3.579244584950415"
$ws.Range("M8").Value = "[-5.673816548749891, 12.832305718650723]"
$ws.Range("N8").Value = [double]"0.4400046642145576"
$ws.Range("O8").Value = [double]"0.4400046642145576"
$ws.Range("P8").Value = [double]"-2.629000458859927"
$ws.Range("Q8").Value = "[-5.767448375058739, 0.5094474573388856]"
$ws.Range("R8").Value = [double]"0.09849279007897094"
$ws.Range("S8").Value = [double]"0.09849279007897094"
$ws.Range("T8").Value = [double]"14.89721032671989"
$ws.Range("U8").Value = "[9.947708992730291, 19.84671166070948]"
$ws.Range("V8").Value = [double]"2.524813120796665e-07"
$ws.Range("W8").Value = [double]"2.524813120796665e-07"
$ws.Range("X8").Value = [double]"10.76590590590615"
$ws.Range("Y8").Value = [double]"-2.086216216216265"
$ws.Range("Z8").Value = [double]"23.61802802802856"
$ws.Range("B9").Value = [double]"0"
$ws.Range("F9").Value = [double]"24.38000000000037"
$ws.Range("H9").Value = [double]"0.1972419342483075"
$ws.Range("I9").Value = [double]"0.1972419342483075"
$ws.Range("L9").Value = [double]"6.263827159303062"
$ws.Range("M9").Value = "[-2.870493690099666, 15.39814800870579]"
$ws.Range("N9").Value = [double]"0.1740490337142737"
$ws.Range("O9").Value = [double]"0.1740490337142737"
$ws.Range("P9").Value = [double]"2.484342538874811"
$ws.Range("Q9").Value = "[-0.6163685286322318, 5.585053606381854]"
$ws.Range("R9").Value = [double]"0.1135766847836717"
$ws.Range("S9").Value = [double]"0.1135766847836717"
$ws.Range("T9").Value = [double]"14.57507059333473"
$ws.Range("U9").Value = "[9.717698758776613, 19.432442427892845]"
$ws.Range("V9").Value = [double]"2.690581957232041e-07"
$ws.Range("W9").Value = [double]"2.690581957232041e-07"
$ws.Range("X9").Value = [double]"14.74026026026049"
$ws.Range("Y9").Value = [double]"2.708888888888932"
$ws.Range("Z9").Value = [double]"26.77163163163204"
$ws.Range("F10").Value = [double]"24.38000000000037"
$ws.Range("H10").Value = [double]"0.1033446231017542"
$ws.Range("I10").Value = [double]"0.1033446231017542"
$ws.Range("L10").Value = [double]"8.630167181353993"
$ws.Range("M10").Value = "[-1.058040527195697, 18.318374889903684]"
$ws.Range("N10").Value = [double]"0.07951016790290244"
$ws.Range("O10").Value = [double]"0.07951016790290244"
$ws.Range("P10").Value = [double]"1.591237119836271"
$ws.Range("Q10").Value = "[0.11950002085726918, 3.062974218815273]"
$ws.Range("R10").Value = [double]"0.03471500659321003"
$ws.Range("S10").Value = [double]"0.03471500659321003"
$ws.Range("T10").Value = [double]"13.84651548690181"
$ws.Range("U10").Value = "[8.180437463578844, 19.512593510224768]"
$ws.Range("V10").Value = [double]"1.19226894756963e-05"
$ws.Range("W10").Value = [double]"1.19226894756963e-05"
$ws.Range("X10").Value = [double]"18.20568568568596"
$ws.Range("Y10").Value = [double]"12.49505505505525"
$ws.Range("Z10").Value = [double]"23.91631631631668"
$ws.Range("F11").Value = [double]"24.38000000000037"
$ws.Range("H11").Value = [double]"0.02499182444089931"
$ws.Range("I11").Value = [double]"0.02499182444089931"
$ws.Range("L11").Value = [double]"7.068916448899532"
$ws.Range("M11").Value = "[0.24011046335190045, 13.897722434447163]"
$ws.Range("N11").Value = [double]"0.04278017344429741"
$ws.Range("O11").Value = [double]"0.04278017344429741"
$ws.Range("P11").Value = [double]"2.345974093671657"
$ws.Range("Q11").Value = "[1.50318447288881, 3.1887637144545034]"
$ws.Range("R11").Value = [double]"1.194219149347475e-06"
$ws.Range("S11").Value = [double]"1.194219149347475e-06"
$ws.Range("T11").Value = [double]"13.70166544819072"
$ws.Range("U11").Value = "[10.14132035150768, 17.26201054487376]"
$ws.Range("V11").Value = [double]"7.957399184022051e-10"
$ws.Range("W11").Value = [double]"7.957399184022051e-10"
$ws.Range("X11").Value = [double]"15.27715715715739"
$ws.Range("Y11").Value = [double]"12.00696696696716"
$ws.Range("Z11").Value = [double]"18.54734734734763"
$ws.Range("F12").Value = [double]"24.38000000000037"
$ws.Range("H12").Value = [double]"0.1032186183284136"
$ws.Range("I12").Value = [double]"0.1032186183284136"
$ws.Range("L12").Value = [double]"8.347767463108028"
$ws.Range("M12").Value = "[-1.9457161306254775, 18.641251056841533]"
$ws.Range("N12").Value = [double]"0.1093654711467971"
$ws.Range("O12").Value = [double]"0.1093654711467971"
$ws.Range("P12").Value = [double]"2.182447749340657"
$ws.Range("Q12").Value = "[-0.6918422260157708, 5.056737724697085]"
$ws.Range("R12").Value = [double]"0.1331877998710298"
$ws.Range("S12").Value = [double]"0.1331877998710298"
$ws.Range("T12").Value = [double]"15.77058854573799"
$ws.Range("U12").Value = "[10.309193296475094, 21.23198379500089]"
$ws.Range("V12").Value = [double]"5.851655451305504e-07"
$ws.Range("W12").Value = [double]"5.851655451305504e-07"
$ws.Range("X12").Value = [double]"15.91167167167192"
$ws.Range("Y12").Value = [double]"4.758858858858932"
$ws.Range("Z12").Value = [double]"27.0644844844849"
$ws.Range("F13").Value = [double]"24.38000000000037"
$ws.Range("H13").Value = [double]"0.4465097247377168"
$ws.Range("I13").Value = [double]"0.4465097247377168"
$ws.Range("L13").Value = [double]"5.001766328331279"
$ws.Range("M13").Value = "[-5.165873697323401, 15.16940635398596]"
$ws.Range("N13").Value = [double]"0.3270822606842252"
$ws.Range("O13").Value = [double]"0.3270822606842252"
$ws.Range("P13").Value = [double]"1.855395060678656"
$ws.Range("Q13").Value = "[-1.2704739059562327, 4.981264027313546]"
$ws.Range("R13").Value = [double]"0.2381562100039811"
$ws.Range("S13").Value = [double]"0.2381562100039811"
$ws.Range("T13").Value = [double]"16.88187112750275"
$ws.Range("U13").Value = "[11.287199661315672, 22.47654259368983]"
$ws.Range("V13").Value = [double]"2.395193996473921e-07"
$ws.Range("W13").Value = [double]"2.395193996473921e-07"
$ws.Range("X13").Value = [double]"17.18070070070096"
$ws.Range("Y13").Value = [double]"5.051711711711789"
$ws.Range("Z13").Value = [double]"29.30968968969014"
$ws.Range("F14").Value = [double]"24.38000000000037"
$ws.Range("H14").Value = [double]"0.5141420459001139"
$ws.Range("I14").Value = [double]"0.5141420459001139"
$ws.Range("L14").Value = [double]"4.1073042727351"
$ws.Range("M14").Value = "[-5.560723685520521, 13.77533223099072]"
$ws.Range("N14").Value = [double]"0.3967196131218853"
$ws.Range("O14").Value = [double]"0.3967196131218853"
$ws.Range("P14").Value = [double]"2.056658253701427"
$ws.Range("Q14").Value = "[-1.0755001877154236, 5.188816695118277]"
$ws.Range("R14").Value = [double]"0.192679506253933"
$ws.Range("S14").Value = [double]"0.192679506253933"
$ws.Range("T14").Value = [double]"12.51795350619663"
$ws.Range("U14").Value = "[7.457094372430916, 17.578812639962337]"
$ws.Range("V14").Value = [double]"9.77188500228543e-06"
$ws.Range("W14").Value = [double]"9.77188500228543e-06"
$ws.Range("X14").Value = [double]"16.39975975976001"
$ws.Range("Y14").Value = [double]"4.246366366366432"
$ws.Range("Z14").Value = [double]"28.55315315315359"

Write-Host "Applied 235 cell updates"
